$d = $word.ActiveDocument

# Locate the paragraph whose text is the "Features" heading, then insert a
# new paragraph right after it (pushing the existing trailing blank
# paragraph further down) introducing the feature list.
$featuresPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd("`r`a") -match "^Features\s*$") {
        $featuresPara = $p
        break
    }
}

# InsertParagraphAfter clones the paragraph/run formatting of $featuresPara
# (jc=both, bold, sz 24) for the freshly inserted paragraph.
$featuresPara.Range.InsertParagraphAfter()

$newPara = $d.Paragraphs.Item($featuresPara.Index + 1)
$newPara.Range.Text = "Following are the features"

Write-Output "inserted features intro paragraph"
